# Update the two-digit multiplication practice sheet to the next day's
# date and a freshly generated set of multiplication problems.

$d = $word.ActiveDocument

# Mapping of old text -> new text (find is literal, not a wildcard match)
$replacements = @(
    @("2023-12-23 Saturday", "2023-12-24 Sunday"),
    @("18×79=", "24×81="),
    @("46×44=", "73×19="),
    @("72×23=", "28×51="),
    @("69×78=", "30×21="),
    @("73×11=", "85×75="),
    @("99×13=", "20×77="),
    @("86×24=", "32×76="),
    @("79×66=", "59×60="),
    @("97×39=", "17×45="),
    @("74×24=", "52×64="),
    @("31×19=", "93×46="),
    @("13×58=", "20×39="),
    @("18×95=", "68×73="),
    @("76×73=", "63×43="),
    @("63×57=", "46×43="),
    @("83×42=", "58×66="),
    @("66×70=", "40×20="),
    @("97×44=", "61×97="),
    @("81×17=", "62×46="),
    @("70×71=", "51×98="),
    @("30×81=", "30×77="),
    @("36×47=", "82×22="),
    @("50×23=", "16×94="),
    @("72×28=", "83×52="),
    @("94×93=", "93×61=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, `
                                  $true, 1, $false, $new, 2)

    if (-not $found) {
        Write-Host "WARNING: could not find text '$old'"
    }
}

$d.Save()
